$p = $ppt.ActivePresentation

# --- Slide 1 title: merge the run-per-word runs into a single run ------
# The shape's current text ("First slide") already equals the desired
# final text, so a direct re-assignment is a textual no-op and the
# multiple <a:r> runs would be left untouched. Route through a
# throwaway value first so the engine is forced to rewrite the
# paragraph into a single run (preserving <a:pPr>/<a:rPr>), then set
# the real text (now trivially a single-run -> single-run change).
$slide1Title = $p.Slides.Item(1).Shapes.Item(1)
$slide1Title.TextFrame.TextRange.Text = "*"
$slide1Title.TextFrame.TextRange.Text = "First slide"

# --- Slide 3 title: same situation ("Third slide") --------------------
$slide3Title = $p.Slides.Item(3).Shapes.Item(1)
$slide3Title.TextFrame.TextRange.Text = "*"
$slide3Title.TextFrame.TextRange.Text = "Third slide"

# --- Speaker notes attached to slide 2: merge the word-per-run notes --
# text into a single run.
$notesBody = $p.Slides.Item(2).NotesPage.Shapes.Item(2)
$notesBody.TextFrame.TextRange.Text = "Some notes here: this first slide should use the Blank template"
